$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: fill in the single transaction/item row.
# Make the "name", "balance" and "transactions count" cells text-formatted
# (matches numFmtId 49 = "@" applied to the underlying style) before writing
# their string values so things like "0:0"/"1:0" are not reinterpreted.
$ws.Range("B4").NumberFormat = "@"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("F4").NumberFormat = "@"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("N4").NumberFormat = "@"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("I4").NumberFormat = "@"
$ws.Range("J4").NumberFormat = "@"
$ws.Range("K4").NumberFormat = "@"

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "BEBY RELIEF 25 MG  SUPP"
$ws.Range("H4").Value = "0:0"
$ws.Range("L4").Value = 38
$ws.Range("N4").Value = "1:0"

# Row 5 (totals row): total price + slightly taller row to match new content.
$ws.Range("K5").Value = 38
$ws.Range("A5").RowHeight = 26.25
